# "Corrected incremental heat rate"
# The unit_incremental_heat_rate operating-point values in column B (rows 6-8)
# were hard-coded numbers; replace them with the correct formulas that derive
# the heat rate as the reciprocal of the (efficiency) value they were
# mistakenly storing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Formula = "=1/0.75"
$ws.Range("B7").Formula = "=1/0.7"
$ws.Range("B8").Formula = "=1/0.667"

# Leave the selection where the author left it when saving.
$ws.Range("C14").Select()
